$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B13: "Disk Clone" -> "Disk to Image Clone"
$ws.Range("B13").Value = "Disk to Image Clone"

# Add new row 16: progress_cloning_drive / Cloning {0} to {1}
$ws.Range("A16").Value = "progress_cloning_drive"
$ws.Range("B16").Value = "Cloning {0} to {1}"

# Update selection to B16
$ws.Range("B16").Select()
